$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price / 1h-volume / hour values scraped on
# Tue Jan 31 03:00:32 UTC 2023. Each target cell already stores its value
# as text (General-formatted numeric-looking strings), so set the number
# format to Text ("@") before writing, exactly as Excel requires to keep
# a numeric-looking literal stored as text instead of auto-converting it.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '309.26'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '-2.00%'
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '38.14'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '-3.34%'
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.061'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-1.37%'
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.07786'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '-4.70%'
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.354'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '-0.53%'
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.896'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '-3.03%'
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '8.183'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '-1.81%'
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '-5.73%'
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '-1.77%'
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.1264'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '-2.89%'
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.1884'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '-4.35%'
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.08804'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-2.73%'
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.09714'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '-0.29%'
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.001374'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '-3.11%'
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '-0.02%'
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.594'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '-1.48%'
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.3411'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '-2.29%'
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.029'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '1.34%'
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '-3.82%'
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.2590'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '3.80%'
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.02104'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '5,587.60%'
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.04382'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '0.49%'
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-2.83%'
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.004263'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-10.45%'
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0001350'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '-65.33%'
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.02141'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '-3.26%'
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.04998'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '-3.24%'
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.008083'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '4.18%'
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.01004'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-3.43%'
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1344'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '-4.05%'
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.002061'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '0.86%'
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.008707'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '-6.18%'
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006433'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '-7.53%'
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000750'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '-0.12%'
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.003439'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '19.23%'
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-0.22%'
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.00002100'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '-0.12%'
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '3'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002000'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '-0.12%'
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '3'

Write-Host "Applied 120 cell updates"
